$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '57.448.85'
$ws.Range("E2").Value2 = '  +1.77%  '

$ws.Range("D3").Value2 = '2.323.30'
$ws.Range("E3").Value2 = '  +0.47%  '

$ws.Range("E4").Value2 = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '544.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = '  +6.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = '  +2.65%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.994'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = '  -0.57%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.538'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = '  +1.03%  '

$ws.Range("D9").Value2 = '2.347.51'
$ws.Range("E9").Value2 = '  +1.42%  '

$ws.Range("E10").Value2 = '  +2.05%  '

$ws.Range("E11").Value2 = '  +0.89%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.42'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = '  +3.66%  '

$ws.Range("E13").Value2 = '  +7.26%  '

$ws.Range("B14").Value2 = 'Avalanche'
$ws.Range("C14").Value2 = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.63'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = '  +1.08%  '

$ws.Range("B15").Value2 = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value2 = '2.741.51'
$ws.Range("E15").Value2 = '  +0.55%  '

$ws.Range("D16").Value2 = '57.531.57'
$ws.Range("E16").Value2 = '  +1.99%  '

$ws.Range("E17").Value2 = '  +1.28%  '

$ws.Range("D18").Value2 = '2.350.66'
$ws.Range("E18").Value2 = '  +1.20%  '

$ws.Range("E19").Value2 = '  +2.57%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '333.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = '  +2.39%  '

$ws.Range("E21").Value2 = '  +2.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = '  +0.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.997'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = '  -0.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = '  +0.39%  '

$ws.Range("E25").Value2 = '  +3.52%  '

$ws.Range("E26").Value2 = '  -0.97%  '

$ws.Range("E27").Value2 = '  -0.32%  '

$ws.Range("E28").Value2 = '  +8.16%  '

$ws.Range("E29").Value2 = '  +5.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.59'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = '  +1.40%  '

$ws.Range("E31").Value2 = '  +2.53%  '

$ws.Range("E32").Value2 = '  +1.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = '  +17.89%  '

$ws.Range("E34").Value2 = '  +1.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = '  -0.04%  '

$ws.Range("B36").Value2 = 'NEARProtocol'
$ws.Range("C36").Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.19'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = '  +7.31%  '

$ws.Range("B37").Value2 = 'FirstDigitalUSD'
$ws.Range("C37").Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.991'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = '  -0.65%  '

$ws.Range("E38").Value2 = '  +2.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.64'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = '  +5.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '39.18'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = '  +1.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '149.05'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = '  -0.23%  '

$ws.Range("E42").Value2 = '  +2.19%  '

$ws.Range("E43").Value2 = '  +1.53%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '284.48'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = '  +3.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.18'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = '  +6.12%  '

$ws.Range("B46").Value2 = 'Hedera'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0506'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = '  +2.63%  '

$ws.Range("B47").Value2 = 'Stellar'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0925'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = '  +0.01%  '

$ws.Range("E48").Value2 = '  +1.74%  '

$ws.Range("B49").Value2 = 'Polygon'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.385'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = '  +9.40%  '

$ws.Range("B50").Value2 = 'EnergySwap'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.63'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = '  +4.32%  '

$ws.Range("B51").Value2 = 'VeChain'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0217'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = '  +1.51%  '
